# Generate Report for Archive
#
# Two localization entries (768147a8-1a1d-441d-a834-4c4de791356f.md and
# 811ad062-771d-4f3b-bdeb-909cb2cda9aa.md) moved from "Ready for handoff"
# to "In Translation". Update the Status columns on every sheet that
# tracks them: the "Overview" sheet (zh-cn / de-de columns) and the
# per-locale "zh-cn" / "de-de" sheets (Status column).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 3 & 4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: Status column C, rows 3 & 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet: Status column C, rows 3 & 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
